# Insert two new data rows (252 and 253) into the sheet, pushing the
# existing rows 252..375 down to 254..377, and populate the two new rows
# with their own data (same constant columns as every other data row,
# plus the specific date/volume/price values for these two new entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 252 onward down by two rows, creating two blank rows.
$ws.Rows("252:253").Insert()

# Row 252 (new) - "Primera" quality entry
$ws.Cells.Item(252, 1).Value = 3
$ws.Cells.Item(252, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(252, 3).Value = "Coquimbo"
$ws.Cells.Item(252, 4).Value = 44460
$ws.Cells.Item(252, 5).Value = 5
$ws.Cells.Item(252, 6).Value = 100112023
$ws.Cells.Item(252, 7).Value = "Brócoli"
$ws.Cells.Item(252, 8).Value = "Sin especificar"
$ws.Cells.Item(252, 9).Value = "Primera"
$ws.Cells.Item(252, 10).Value = 3000
$ws.Cells.Item(252, 11).Value = 500
$ws.Cells.Item(252, 12).Value = 550
$ws.Cells.Item(252, 13).Value = 527
$ws.Cells.Item(252, 14).Value = "$/unidad"
$ws.Cells.Item(252, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(252, 16).Value = 527
$ws.Cells.Item(252, 17).Value = 1
$ws.Cells.Item(252, 18).Value = "Hortaliza"

# Row 253 (new) - "Segunda" quality entry
$ws.Cells.Item(253, 1).Value = 3
$ws.Cells.Item(253, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(253, 3).Value = "Coquimbo"
$ws.Cells.Item(253, 4).Value = 44460
$ws.Cells.Item(253, 5).Value = 5
$ws.Cells.Item(253, 6).Value = 100112023
$ws.Cells.Item(253, 7).Value = "Brócoli"
$ws.Cells.Item(253, 8).Value = "Sin especificar"
$ws.Cells.Item(253, 9).Value = "Segunda"
$ws.Cells.Item(253, 10).Value = 1350
$ws.Cells.Item(253, 11).Value = 400
$ws.Cells.Item(253, 12).Value = 400
$ws.Cells.Item(253, 13).Value = 400
$ws.Cells.Item(253, 14).Value = "$/unidad"
$ws.Cells.Item(253, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(253, 16).Value = 400
$ws.Cells.Item(253, 17).Value = 1
$ws.Cells.Item(253, 18).Value = "Hortaliza"
